{"js": "// Replace the 25 two-digit multiplication problems in the table with a\n// freshly generated set, preserving each cell's run formatting (font,\n// size, paragraph alignment, etc.) \u2014 only the literal text changes.\n\n// New values, indexed by the table's non-blank rows (0, 4, 9, 14, 19)\n// in row-major order, 5 cells per row (matches the grid in the file).\nconst newValues = [\n  [\"64\u00d749=\", \"25\u00d797=\", \"63\u00d764=\", \"87\u00d717=\", \"52\u00d715=\"],\n  [\"26\u00d780=\", \"46\u00d792=\", \"68\u00d755=\", \"59\u00d741=\", \"37\u00d758=\"],\n  [\"65\u00d720=\", \"60\u00d797=\", \"92\u00d789=\", \"17\u00d726=\", \"79\u00d790=\"],\n  [\"82\u00d785=\", \"90\u00d735=\", \"75\u00d793=\", \"15\u00d793=\", \"14\u00d785=\"],\n  [\"31\u00d786=\", \"63\u00d790=\", \"67\u00d723=\", \"70\u00d711=\", \"59\u00d784=\"],\n];\nconst problemRowIndexes = [0, 4, 9, 14, 19];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < problemRowIndexes.length; r++) {\n  const rowIndex = problemRowIndexes[r];\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(rowIndex, c);\n    const range = cell.body.getRange();\n    range.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit multiplication problems in the table with a\n# freshly generated set, preserving each cell's run formatting (font,\n# size, paragraph alignment, etc.) \u2014 only the literal text changes.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New values, indexed by the table's non-blank rows (1, 5, 10, 15, 20 \u2014\n# Word COM is 1-based) in row-major order, 5 cells per row.\n$newValues = @{\n    1  = @(\"64\u00d749=\", \"25\u00d797=\", \"63\u00d764=\", \"87\u00d717=\", \"52\u00d715=\")\n    5  = @(\"26\u00d780=\", \"46\u00d792=\", \"68\u00d755=\", \"59\u00d741=\", \"37\u00d758=\")\n    10 = @(\"65\u00d720=\", \"60\u00d797=\", \"92\u00d789=\", \"17\u00d726=\", \"79\u00d790=\")\n    15 = @(\"82\u00d785=\", \"90\u00d735=\", \"75\u00d793=\", \"15\u00d793=\", \"14\u00d785=\")\n    20 = @(\"31\u00d786=\", \"63\u00d790=\", \"67\u00d723=\", \"70\u00d711=\", \"59\u00d784=\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $values = $newValues[$row]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
